$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.50%'
$ws.Range("H3").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.76%'
$ws.Range("H4").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.22%'
$ws.Range("H5").Value = '震荡市，未发现操作点：平滑 True，成交量无下跌 True，成交量无突破 True'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.39%'
$ws.Range("H6").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.50%'
$ws.Range("H7").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2.11%'
$ws.Range("H8").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.21%'
$ws.Range("H9").Value = '震荡市，未发现操作点：平滑 True，成交量无下跌 True，成交量无突破 True'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.71%'
$ws.Range("H10").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.11%'
$ws.Range("H11").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.43%'
$ws.Range("H12").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.24%'
$ws.Range("H13").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.05%'
$ws.Range("H14").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.52%'
$ws.Range("H15").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.57%'
$ws.Range("H16").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.84%'
$ws.Range("H17").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.62%'
$ws.Range("H18").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.54%'
$ws.Range("H19").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.73%'
$ws.Range("H20").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.70%'
$ws.Range("H21").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.41%'
$ws.Range("H22").Value = '震荡市，低位吸纳 200.00'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.78%'
$ws.Range("H23").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '3.16%'
$ws.Range("H24").Value = '震荡→弱升，试探性建仓 200.00'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.40%'
$ws.Range("H25").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.03%'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2.88%'
$ws.Range("H27").Value = '震荡市，未发现操作点：平滑 False，成交量无下跌 True，成交量无突破 True'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '0.72%'
$ws.Range("H28").Value = '震荡市，未发现操作点：平滑 True，成交量无下跌 True，成交量无突破 True'

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '-0.21%'
$ws.Range("H29").Value = '无操作'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '2.95%'
$ws.Range("H30").Value = '震荡市，高位减持 -10.00% 仓位'
